$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "The consolidation reported insufficient staff to meet this portion of the requirements. ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
